# TC04 Bento multi-filter workbook update:
# "Other treatment given as part of a CTSU protocol" -> "standard AC followed by a taxane"
# for the chemotherapy_regimen filter used across the 4 Cypher queries (Cases/Samples/Files/Summary),
# plus the associated row-height / selection view-state tweaks made while editing in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$oldText = "Other treatment given as part of a CTSU protocol"
$newText = "standard AC followed by a taxane"

# The four query cells (B2:B4 hold the distinct queries, C2:C4 share one query string)
# all contain the WHERE clause with the chemotherapy_regimen filter that changed.
$cellAddresses = @("B2", "C2", "B3", "C3", "B4", "C4")

foreach ($addr in $cellAddresses) {
    $cell = $ws.Range($addr)
    $text = $cell.Value2
    if ($text -ne $null -and $text.Contains($oldText)) {
        $cell.Value2 = $text.Replace($oldText, $newText)
    }
}

# Row 3 (the CasesTab query) shrank slightly in wrapped height after the edit.
$ws.Rows.Item(3).RowHeight = 406

# Selection moved from D10 to C14 while reviewing the sheet.
$null = $ws.Range("C14").Select()
